$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update row 2 election result figures (PORTO / LOUSADA) as per updated tables
$ws.Range("H2").Value = 536
$ws.Range("I2").Value = 1206
$ws.Range("J2").Value = 5194
$ws.Range("K2").Value = 16
$ws.Range("L2").Value = 1439
$ws.Range("M2").Value = 68
$ws.Range("N2").Value = 886
$ws.Range("O2").Value = 4
$ws.Range("P2").Value = 22
$ws.Range("Q2").Value = 12
$ws.Range("R2").Value = 81
$ws.Range("S2").Value = 526
$ws.Range("T2").Value = 925
$ws.Range("U2").Value = 58
$ws.Range("V2").Value = 8028
$ws.Range("W2").Value = 4
$ws.Range("X2").Value = 8187
$ws.Range("Z2").Value = 133
